{"js": "// Append the 2018-03-18 work-log update to the end of the \"\uc774\ubc88\uc8fc \uc791\uc5c5\uc77c\uc9c0\"\n// paragraph that currently ends with \"...\ub4a4\uc5d0\uc788\ub294\uac83\ub3c4 \uc548\uadf8\ub824\uc9d0.\" in the\n// second weekly table (2018-03-12 ~ 2018-03-18).\n//\n// The appended text is simply tacked on right after the existing\n// \" \uc548\uadf8\ub824\uc9d0.\" run, turning it into:\n//   \" \uc548\uadf8\ub824\uc9d0. \ub9ac\uc9c0\ub4dc\ubc14\ub514 1\ucc28 \ud14c\uc2a4\ud2b8\uc6a9 \ucd94\uac00. \ud1a0\ud06c\ub97c \uc0ac\uc6a9\ud558\ub294 \uac83 \uae4c\uc9c4 \ub41c\ub2e4.\"\n\nconst body = context.document.body;\n\n// \"\uc548\uadf8\ub824\uc9d0.\" is unique in this document, so a plain text search locates\n// the exact run/range we need to extend.\nconst results = body.search(\"\uc548\uadf8\ub824\uc9d0.\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find target text \"\uc548\uadf8\ub824\uc9d0.\" in the document.');\n}\n\nconst target = results.items[0];\ntarget.insertText(\n  \" \ub9ac\uc9c0\ub4dc\ubc14\ub514 1\ucc28 \ud14c\uc2a4\ud2b8\uc6a9 \ucd94\uac00. \ud1a0\ud06c\ub97c \uc0ac\uc6a9\ud558\ub294 \uac83 \uae4c\uc9c4 \ub41c\ub2e4.\",\n  \"After\"\n);\n\nawait context.sync();\n", "ps1": "# Append the 2018-03-18 work-log update right after the existing\n# \" \uc548\uadf8\ub824\uc9d0.\" text in the \"\uc774\ubc88\uc8fc \uc791\uc5c5\uc77c\uc9c0\" cell of the second weekly\n# table (2018-03-12 ~ 2018-03-18), turning it into:\n#   \" \uc548\uadf8\ub824\uc9d0. \ub9ac\uc9c0\ub4dc\ubc14\ub514 1\ucc28 \ud14c\uc2a4\ud2b8\uc6a9 \ucd94\uac00. \ud1a0\ud06c\ub97c \uc0ac\uc6a9\ud558\ub294 \uac83 \uae4c\uc9c4 \ub41c\ub2e4.\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"\uc548\uadf8\ub824\uc9d0.\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n\nif ($find.Execute()) {\n    $rng = $find.Parent\n    $rng.Collapse(0)  # wdCollapseEnd\n    $rng.InsertAfter(\" \ub9ac\uc9c0\ub4dc\ubc14\ub514 1\ucc28 \ud14c\uc2a4\ud2b8\uc6a9 \ucd94\uac00. \ud1a0\ud06c\ub97c \uc0ac\uc6a9\ud558\ub294 \uac83 \uae4c\uc9c4 \ub41c\ub2e4.\")\n}\n"}
